$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Run 2" (row 3) result values (previously only G3/H3 had
#     preliminary numbers, now all four metrics are final) ---
$ws.Range("E3").Value = 0.364595036529411
$ws.Range("F3").Value = 0.403574780101712
$ws.Range("G3").Value = 0.509321365387381
$ws.Range("H3").Value = 0.579446402898139

# --- Re-format the results block to 4 decimal places, centred.
#     Done per contiguous block (rather than one comma-joined range) so the
#     formatting reliably reaches every cell. ---
$ws.Range("E2:H3").HorizontalAlignment = -4108
$ws.Range("E2:H3").NumberFormat = "0.0000"

$ws.Range("E4:F4").HorizontalAlignment = -4108
$ws.Range("E4:F4").NumberFormat = "0.0000"

# G4:H4 additionally keeps its wrapped-text styling on top of the new format
$ws.Range("G4:H4").HorizontalAlignment = -4108
$ws.Range("G4:H4").NumberFormat = "0.0000"
$ws.Range("G4:H4").WrapText = $true

# --- Reserve the blank rows (5-18) for the runs still to come, formatted
#     the same way as the rest of the results block ---
$ws.Range("E5:H18").NumberFormat = "0.0000"
$ws.Range("E5:H18").HorizontalAlignment = -4108

# --- Row heights for the data rows ---
$ws.Range("A2:A18").EntireRow.RowHeight = 13.8

# --- Sheet view bookkeeping to match the saved state ---
$ws.Range("B15").Select() | Out-Null
